$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.330.01'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.358.37'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.65'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.55'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').Value = '  +0.61%  '
$ws.Range('E8').Value = '  +4.82%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  +5.70%  '
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  +2.78%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.85'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.775.51'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.285.49'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.381.07'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.73'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '333.12'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.27'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.78'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '62.90'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.07%  '
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('E25').Value = '  -2.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('E27').Value = '  +5.38%  '
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.99'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.14'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('E33').Value = '  +12.51%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.27'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +6.95%  '
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('E38').Value = '  +4.24%  '
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '143.45'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '293.86'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.379'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('E43').Value = '  +1.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0949'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.21'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('E47').Value = '  +1.56%  '
$ws.Range('E48').Value = '  +1.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.386'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.51'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('E51').Value = '  +0.52%  '
